# Applies the "added marks for 26-05-2025" commit:
#  1. Fills in the previously-blank Physics/Biology/Chemistry rows (56-58,
#     corresponding to date 26-05-2025 / serial 45803) with their
#     correct/incorrect/unattempted counts.
#  2. Corrects F2 (incorrect count for the first row) from 14 to 15.
#  3. Adds the missing `no_of_questions` calculated-column formula
#     (=SUM(Table1[[#This Row],[correct]:[unattempted]])) to column D for
#     every data row of Table1 (rows 2-58), matching the rest of the
#     calculated table columns.
#  4. Updates the sheet view (scroll position / selection) to reflect
#     where the author was last working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Marksheet")

# ---------------------------------------------------------------------
# 1. New data for 26-05-2025 (rows 56-58: Physics, Biology, Chemistry)
# ---------------------------------------------------------------------
$newData = @{
  56 = @(47, 2, 1)   # Physics:   correct, incorrect, unattempted
  57 = @(19, 2, 0)   # Biology:   correct, incorrect, unattempted
  58 = @(47, 3, 0)   # Chemistry: correct, incorrect, unattempted
}

foreach ($r in $newData.Keys) {
  $vals = $newData[$r]
  $ws.Range("E$r").Value = $vals[0]
  $ws.Range("F$r").Value = $vals[1]
  $ws.Range("G$r").Value = $vals[2]
}

# ---------------------------------------------------------------------
# 2. Correct the incorrect-answer count for row 2 (14 -> 15)
# ---------------------------------------------------------------------
$ws.Range("F2").Value = 15

# ---------------------------------------------------------------------
# 3. Add the calculated `no_of_questions` formula to column D for every
#    row of the table (2-58). D5 currently has no explicit style (style
#    0) while every sibling cell uses style index 5 (thin border all
#    around); copy the formatting from a neighboring cell first so the
#    style matches after the formula is applied.
# ---------------------------------------------------------------------
$formula = "=SUM(Table1[[#This Row],[correct]:[unattempted]])"

$ws.Range("E5").Copy()
$ws.Range("D5").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

for ($r = 2; $r -le 58; $r++) {
  $ws.Range("D$r").Formula = $formula
}

# ---------------------------------------------------------------------
# Force a clean recalculation of the dependent calculated columns for
# the newly populated rows (56-58). These cells previously evaluated to
# errors (#DIV/0!, #VALUE!) while D/E/F/G were blank, and the cached
# error-result type needs to be explicitly refreshed now that valid
# numbers are present.
# ---------------------------------------------------------------------
foreach ($r in @(56, 57, 58)) {
  foreach ($col in @("H", "I", "J", "K", "L", "M", "N")) {
    $addr = "$col$r"
    $f = $ws.Range($addr).Formula
    $ws.Range($addr).Formula = $f
  }
}

$excel.CalculateFullRebuild()

# ---------------------------------------------------------------------
# 4. Update sheet view: scroll position and active selection
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 51
$win.ScrollColumn = 1
$ws.Range("D57").Select()
